$wb = $excel.ActiveWorkbook

# --- During the presentation, F11 was pressed (new Chart sheet with the
#     default chart type) while a single cell was selected on Tabelle1.
#     This inserts a brand-new chart sheet "Diagramm1" as the very first
#     sheet in the workbook. ---
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, [System.Type]::Missing, [System.Type]::Missing, -4109)
$newSheet.Name = "Diagramm1"

# The default chart Excel inserts via F11 (no data selected) is an empty
# clustered column chart.
$shp = $newSheet.Shapes.AddChart2(201, 51)
$shp.Chart.ChartType = 51

# --- Back to Tabelle1: re-select it (so it becomes the active tab again)
#     and leave the cell selection where the presenter left it. ---
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws1.Activate()
$ws1.Range("J4").Select()
